$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.213374376296997
$ws.Range("B1").Value = 2.391945600509644
$ws.Range("D1").Value = 1.417901754379272
$ws.Range("E1").Value = 0.9105518460273743
